# Reverse the order of the "Periodo Mora" / "Valor Mora" table (rows 16-40)
# so the most recent period (2408) appears first and the oldest (2208) last.
# This mirrors the real edit: the underlying period list was rebuilt in
# descending order, carrying the one differing "Valor Mora" (36400, tied to
# period 2408) along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("2208","2209","2210","2211","2212","2301","2302","2303","2304","2305","2306","2307","2308","2309","2310","2311","2312","2401","2402","2403","2404","2405","2406","2407","2408")
$values  = @(52000,52000,52000,52000,52000,52000,52000,52000,52000,52000,52000,52000,52000,52000,52000,52000,52000,52000,52000,52000,52000,52000,52000,52000,36400)

$startRow = 16
$count = $periods.Length

for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i
    $srcIdx = $count - 1 - $i
    $ws.Range("E$row").Value = $periods[$srcIdx]
    $ws.Range("F$row").Value = $values[$srcIdx]
}
